$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be stored as Text so that numeric-looking
# strings (e.g. "11.20", "0.600", "58.958.85") keep their exact original
# formatting instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.091.96"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.533.06"
$ws.Range("E3").Value = "  +2.97%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "540.77"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "143.84"
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "2.529.88"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "5.55"
$ws.Range("E12").Value = "  +4.64%  "
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "2.950.25"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").Value = "23.59"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").Value = "59.034.47"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "2.522.08"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "11.20"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "4.28"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").Value = "324.61"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +3.16%  "
$ws.Range("D23").Value = "5.77"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "61.87"
$ws.Range("D25").Value = "0.439"
$ws.Range("E25").Value = "  -4.95%  "
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").Value = "7.90"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").Value = [string]::Concat("0.0", [char]8323, "0777")
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "1.81"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").Value = "6.66"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  -9.70%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "1.47"
$ws.Range("E34").Value = "  +6.46%  "
$ws.Range("D35").Value = "157.23"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").Value = "18.66"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "4.36"
$ws.Range("E37").Value = "  -4.08%  "
$ws.Range("D38").Value = "1.62"
$ws.Range("E38").Value = "  -7.48%  "
$ws.Range("D39").Value = "37.04"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("E40").Value = "  -4.98%  "
$ws.Range("D41").Value = "295.79"
$ws.Range("E41").Value = "  -7.67%  "
$ws.Range("D42").Value = "3.71"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "0.823"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "0.600"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("D46").Value = "10.78"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "0.0929"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "18.64"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "122.53"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").Value = "0.0515"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0228"
$ws.Range("E51").Value = "  -0.97%  "
